# Sync automatico del tracker (cada 3h)
# Appends new rows of match data (rows 153-160) to the tracker sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("14316299", "2025-08-10", "Francisco Comesaña", "Luciano Darderi", "Gana Luciano Darderi", 2.2, "", ""),
    @("14316450", "2025-08-10", "Jasmine Paolini", "Maria Sakkari", "Gana Jasmine Paolini", 1.67, "", ""),
    @("14316458", "2025-08-10", "Veronika Kudermetova", "Belinda Bencic", "Gana Belinda Bencic", 1.62, "", ""),
    @("14316464", "2025-08-10", "Catherine McNally", "McCartney Kessler", "Gana McCartney Kessler", 1.91, "", ""),
    @("14393241", "2025-08-11", "Nicolai Budkov Kjaer", "Juan Bautista Torres", "Gana Juan Bautista Torres", 4.5, "", ""),
    @("14392622", "2025-08-11", "Dan Added", "Aristotelis Thanos", "Gana Dan Added", 2.1, "", ""),
    @("14392524", "2025-08-11", "Jelle Sels", "Daniel Rincon", "Gana Jelle Sels", 2.63, "", ""),
    @("14392521", "2025-08-11", "Luciano Emanuel Ambrogi", "Joel Schwaerzler", "Gana Luciano Emanuel Ambrogi", 5, "", "")
)

$startRow = 153
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    # Column A (event_id) is stored as text in the new rows.
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $data[0]

    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
    $ws.Cells.Item($r, 8).Value = $data[7]
}
